# Update "想去人数" (wanted-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect the latest scrape (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F4").Value = 6201
    $ws.Range("F11").Value = 951
    $ws.Range("F12").Value = 213
    $ws.Range("F13").Value = 5571
}
